# C5-PowerPoint.pptx edit replay
# 1) Re-style the single table on slide 6 with the built-in table style
#    {3D6D102F-D3B1-442C-8DA8-0B8EC3BBAD8A} (was the custom "Table_0" style
#    {A74C3168-DAB0-477E-BE15-554AAA2B6390} defined in tableStyles.xml).
# 2) Swap the deck's theme colours from the custom "Integral" palette over
#    to the stock Office palette (the paired notes-master theme holding the
#    old Office palette is reassigned as a consequence of the same design
#    change on the real file; this host only exposes the single shared
#    ThemeColorScheme, so we drive that).

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{3D6D102F-D3B1-442C-8DA8-0B8EC3BBAD8A}")

# --- 2. Theme colours -------------------------------------------------------
$themeColors = $p.Slides.Item(1).ThemeColorScheme

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -> stock "Office" palette
# (000000, FFFFFF, 44546A, E7E6E6, 5B9BD5, ED7D31, A5A5A5, FFC000, 4472C4,
#  70AD47, 0563C1, 954F72) expressed as BGR-packed RGB() integers.
$officePalette = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $officePalette.Length; $i++) {
    $themeColors.Colors($i).RGB = $officePalette[$i - 1]
}
